$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SWR")

# --- Rename recipe identifiers to uppercase (column A, rows 2-6) ---
$ws.Range("A2").Value = "123456-EXTRAPLACE"
$ws.Range("A3").Value = "123457-NOPLACE-ALL"
$ws.Range("A4").Value = "123458-NOPLACE-PARTIAL"
$ws.Range("A5").Value = "123459-PARTSUB-ALL"
$ws.Range("A6").Value = "123460-PARTSUB-PARTIAL"

# These cells previously carried the "@" text style (s="1"); the
# identifiers no longer need a forced text format, so drop back to
# the sheet's default (General) number format.
$ws.Range("A2:A6").NumberFormat = "General"

# --- Add helper formulas that upper-case the new SWR rows (8-11) ---
$ws.Range("C18").Formula = "=UPPER(A8)"
$ws.Range("C19").Formula = "=UPPER(A9)"
$ws.Range("C20").Formula = "=UPPER(A10)"
$ws.Range("C21").Formula = "=UPPER(A11)"

# --- The workbook no longer needs the deep row-outline grouping ---
$ws.Cells.ClearOutline()
$wb.Worksheets.Item("settings").Cells.ClearOutline()

# --- Duplicate-check formatting should skip the (now-static) header
#     rows 2-6 and resume from row 7 onward, keeping row 1 covered ---
$origCf = $ws.Cells.FormatConditions.Item(1)
$origCf.ModifyAppliesToRange($ws.Range("A7:A1048576"))
$extraCf = $ws.Range("A1").FormatConditions.AddUniqueValues()
$extraCf.DupeUnique = 1
$extraCf.Font.Color = 393372
$extraCf.Interior.Color = 13551615

# --- Font normalisation: use Calibri (Latin) instead of the cached
#     SimSun/宋体 name throughout both sheets ---
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Font.Name = "Calibri"
}

# --- Cosmetic: rename the built-in (localized) cell styles back to
#     their English display names ---
$builtinStyleNames = @(
    "Normal","Comma","Currency","Percent","Comma [0]","Currency [0]",
    "Hyperlink","Followed Hyperlink","Note","Warning Text","Title",
    "CExplanatory Text","Heading 1","Heading 2","Heading 3","Heading 4",
    "Input","Output","Calculation","Check Cell","Linked Cell","Total",
    "Good","Bad","Neutral","Accent1","20% - Accent1","40% - Accent1",
    "60% - Accent1","Accent2","20% - Accent2","40% - Accent2","60% - Accent2",
    "Accent3","20% - Accent3","40% - Accent3","60% - Accent3","Accent4",
    "20% - Accent4","40% - Accent4","60% - Accent4","Accent5","20% - Accent5",
    "40% - Accent5","60% - Accent5","Accent6","20% - Accent6","40% - Accent6",
    "60% - Accent6"
)
for ($i = 1; $i -le $wb.Styles.Count; $i++) {
    try {
        $wb.Styles.Item($i).Name = $builtinStyleNames[$i - 1]
    } catch {
    }
}

# --- Window / selection bookkeeping ---
[void]$ws.Range("F13").Select()
try { $excel.ActiveWindow.Height = 9144 } catch { }
